# Update TPM-derived values in the LR-pairs worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Myoc/Fzd10 -> MuSCs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1603853333333333
$ws.Range("H2").Value = 0.481156
$ws.Range("I2").Value = 0.01032935781992836
$ws.Range("J2").Value = 0.01042870175281933
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.016657
$ws.Range("N2").Value = 0.033314
$ws.Range("Q2").Value = 0.002671538497333334
$ws.Range("R2").Value = 0.016029230984
$ws.Range("S2").Value = 0.01032935781992836
$ws.Range("T2").Value = 0.01042870175281933

# Row 3 (FAPs -> Myoc/Fzd10 -> MuSCs)
$ws.Range("I3").Value = 0.9610926076617912
$ws.Range("J3").Value = 0.9703360399430661
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.016657
$ws.Range("N3").Value = 0.033314
$ws.Range("Q3").Value = 0.2485726553026667
$ws.Range("R3").Value = 1.491435931816
$ws.Range("S3").Value = 0.9610926076617912
$ws.Range("T3").Value = 0.9703360399430661

# Row 4 (MuSCs -> Myoc/Fzd10 -> MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.443735
$ws.Range("H4").Value = 0.88747
$ws.Range("I4").Value = 0.02857803451828042
$ws.Range("J4").Value = 0.01923525830411462
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.016657
$ws.Range("N4").Value = 0.033314
$ws.Range("Q4").Value = 0.007391293895
$ws.Range("R4").Value = 0.02956517558
$ws.Range("S4").Value = 0.02857803451828042
$ws.Range("T4").Value = 0.01923525830411462
